# Update the "PIVL-TS" StructureDefinition spreadsheet:
#  - bump Version and Date on the Metadata sheet
#  - insert a new "Jurisdiction" property row right after "Contact"

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# Insert a new row after "Contact" (row 10) to hold the "Jurisdiction" property.
# This shifts every row below it (Description, Purpose, Copyright, ...) down by one.
$ws.Rows.Item(11).Insert()

# Match the look (borders/shading/alignment) of the other property rows by
# copying the formatting from the row directly above (Contact).
$ws.Range("A10:B10").Copy()
$ws.Range("A11:B11").PasteSpecial(-4122)

$ws.Range("A11").Value = "Jurisdiction"
$ws.Range("B11").Value = ""

# Bump the Version property.
$ws.Range("B3").Value = "2.0.1-sd-202510-matchbox-patch"

# Bump the Date property.
$ws.Range("B8").Value = "2025-10-29T22:15:57+01:00"
